$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" ------------------
# Overview sheet keeps one status column per locale (E = zh-cn, F = de-de).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# Per-locale detail sheets keep the status in column C.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the status columns ----------------------------------------------
# Target stored width ~13.41 chars (was ~17.22 chars). ColumnWidth is quoted
# in "characters" and gets snapped to the workbook's pixel grid on save, so
# feed it the character width whose rounded pixel grid lands closest to the
# desired stored width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
